$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5108
$ws.Range("J40").Value = 6694.25
$ws.Range("L40").Value = 6694.25
$ws.Range("N40").Value = -7044.25
$ws.Range("H51").Value = 17374.875
$ws.Range("I51").Value = 17999.857
$ws.Range("K51").Value = 17999.857
$ws.Range("M51").Value = -17515.857
$ws.Range("H74").Value = 20360.666
$ws.Range("J74").Value = 34000
$ws.Range("L74").Value = 34000
$ws.Range("N74").Value = -35872
$ws.Range("H77").Value = 20360.666
$ws.Range("J77").Value = 34000
$ws.Range("L77").Value = 170000
$ws.Range("N77").Value = -179360
$ws.Range("H106").Value = 7754.769
$ws.Range("I106").Value = 4619.1816
$ws.Range("K106").Value = 4619.1816
$ws.Range("M106").Value = -3988.1816
$ws.Range("H113").Value = 8537
$ws.Range("I113").Value = 9799.571
$ws.Range("J113").Value = 6769.4
$ws.Range("K113").Value = 9799.571
$ws.Range("L113").Value = 6769.4
$ws.Range("M113").Value = -6545.571
$ws.Range("N113").Value = -13277.4
$ws.Range("H116").Value = 18281.834
$ws.Range("I116").Value = 15897.2
$ws.Range("K116").Value = 15897.2
$ws.Range("M116").Value = -12455.2
$ws.Range("H132").Value = 6120.6665
$ws.Range("J132").Value = 2650
$ws.Range("L132").Value = 7950
$ws.Range("N132").Value = -13010
$ws.Range("H135").Value = 892.44446
$ws.Range("I135").Value = 339
$ws.Range("J135").Value = 1999.3334
$ws.Range("K135").Value = 3051
$ws.Range("L135").Value = 17994.0006
$ws.Range("M135").Value = -516
$ws.Range("N135").Value = -23064.0006
$ws.Range("H137").Value = 2531.7
$ws.Range("I137").Value = 1563.5
$ws.Range("J137").Value = 3984
$ws.Range("K137").Value = 4690.5
$ws.Range("L137").Value = 11952
$ws.Range("M137").Value = -2140.5
$ws.Range("N137").Value = -17052
$ws.Range("H138").Value = 4777.3335
$ws.Range("I138").Value = 5171.6665
$ws.Range("K138").Value = 15514.9995
$ws.Range("M138").Value = -10374.9995

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3721.1333
$ws.Range("I45").Value = 2096.5
$ws.Range("K45").Value = 2096.5
$ws.Range("M45").Value = -1719.5
$ws.Range("H61").Value = 3280.625
$ws.Range("I61").Value = 3811.25
$ws.Range("J61").Value = 2750
$ws.Range("K61").Value = 3811.25
$ws.Range("L61").Value = 2750
$ws.Range("M61").Value = -3599.25
$ws.Range("N61").Value = -3174
$ws.Range("H74").Value = 2254.5625
$ws.Range("J74").Value = 2999.875
$ws.Range("L74").Value = 2999.875
$ws.Range("N74").Value = -4747.875
$ws.Range("H77").Value = 2254.5625
$ws.Range("J77").Value = 2999.875
$ws.Range("L77").Value = 14999.375
$ws.Range("N77").Value = -23735.375
$ws.Range("H97").Value = 6448.591
$ws.Range("I97").Value = 9166.416999999999
$ws.Range("J97").Value = 3187.2
$ws.Range("K97").Value = 9166.416999999999
$ws.Range("L97").Value = 3187.2
$ws.Range("M97").Value = -8670.416999999999
$ws.Range("N97").Value = -4179.2
$ws.Range("H102").Value = 3437.7
$ws.Range("I102").Value = 3424.111
$ws.Range("J102").Value = 3560
$ws.Range("K102").Value = 3424.111
$ws.Range("L102").Value = 3560
$ws.Range("M102").Value = -1802.111
$ws.Range("N102").Value = -6804
$ws.Range("H132").Value = 93533.55
$ws.Range("I132").Value = 113718.89
$ws.Range("J132").Value = 2699.5
$ws.Range("K132").Value = 341156.67
$ws.Range("L132").Value = 8098.5
$ws.Range("M132").Value = -338626.67
$ws.Range("N132").Value = -13158.5
$ws.Range("H136").Value = 3280.625
$ws.Range("I136").Value = 3811.25
$ws.Range("J136").Value = 2750
$ws.Range("K136").Value = 11433.75
$ws.Range("L136").Value = 8250
$ws.Range("M136").Value = -8883.75
$ws.Range("N136").Value = -13350

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3848449.8
$ws.Range("I20").Value = 7144052.5
$ws.Range("K20").Value = 7144052.5
$ws.Range("M20").Value = -7143805.5
$ws.Range("H50").Value = 74998.5
$ws.Range("J50").Value = 74998.5
$ws.Range("L50").Value = 74998.5
$ws.Range("N50").Value = -76146.5
$ws.Range("H134").Value = 2940.838
$ws.Range("I134").Value = 2061.5173
$ws.Range("J134").Value = 6128.375
$ws.Range("K134").Value = 6184.5519
$ws.Range("L134").Value = 18385.125
$ws.Range("M134").Value = -3649.5519
$ws.Range("N134").Value = -23455.125

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1298.25
$ws.Range("I105").Value = 1298.25
$ws.Range("K105").Value = 1298.25
$ws.Range("M105").Value = 448.75
$ws.Range("H120").Value = 67999.39999999999
$ws.Range("J120").Value = 74999.25
$ws.Range("L120").Value = 74999.25
$ws.Range("N120").Value = -82257.25

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 63962.09
$ws.Range("J37").Value = 63962.09
$ws.Range("L37").Value = 191886.27
$ws.Range("N37").Value = -192110.27

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 69999
$ws.Range("J94").Value = 69999
$ws.Range("L94").Value = 69999
$ws.Range("N94").Value = -71351
$ws.Range("H102").Value = 3692.4443
$ws.Range("J102").Value = 4242.273
$ws.Range("L102").Value = 4242.273
$ws.Range("N102").Value = -7486.273
$ws.Range("H122").Value = 5007.2
$ws.Range("I122").Value = 4179.5
$ws.Range("K122").Value = 12538.5
$ws.Range("M122").Value = -10088.5
$ws.Range("H126").Value = 6978.4
$ws.Range("J126").Value = 6978.4
$ws.Range("L126").Value = 20935.2
$ws.Range("N126").Value = -25875.2
$ws.Range("H132").Value = 103042
$ws.Range("I132").Value = 145486.86
$ws.Range("K132").Value = 436460.58
$ws.Range("M132").Value = -433930.58

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7433.625
$ws.Range("I40").Value = 3249.75
$ws.Range("K40").Value = 3249.75
$ws.Range("M40").Value = -3113.75
$ws.Range("H46").Value = 5178.5117
$ws.Range("I46").Value = 37670
$ws.Range("J46").Value = 2741.65
$ws.Range("K46").Value = 37670
$ws.Range("L46").Value = 2741.65
$ws.Range("M46").Value = -37482
$ws.Range("N46").Value = -3117.65
$ws.Range("H68").Value = 4219.75
$ws.Range("I68").Value = 3041.1428
$ws.Range("J68").Value = 5869.8
$ws.Range("K68").Value = 3041.1428
$ws.Range("L68").Value = 5869.8
$ws.Range("M68").Value = -2292.1428
$ws.Range("N68").Value = -7367.8
$ws.Range("H71").Value = 4219.75
$ws.Range("I71").Value = 3041.1428
$ws.Range("J71").Value = 5869.8
$ws.Range("K71").Value = 15205.714
$ws.Range("L71").Value = 29349
$ws.Range("M71").Value = -11461.714
$ws.Range("N71").Value = -36837
$ws.Range("H100").Value = 4470.88
$ws.Range("I100").Value = 4245.375
$ws.Range("J100").Value = 4871.778
$ws.Range("K100").Value = 4245.375
$ws.Range("L100").Value = 4871.778
$ws.Range("M100").Value = -3704.375
$ws.Range("N100").Value = -5953.778
$ws.Range("H122").Value = 3460.111
$ws.Range("I122").Value = 2786.56
$ws.Range("J122").Value = 4990.909
$ws.Range("K122").Value = 8359.68
$ws.Range("L122").Value = 14972.727
$ws.Range("M122").Value = -5909.68
$ws.Range("N122").Value = -19872.727
$ws.Range("H136").Value = 3699.4285
$ws.Range("I136").Value = 2566
$ws.Range("K136").Value = 7698
$ws.Range("M136").Value = -5148

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 58191.5
$ws.Range("I45").Value = 70000
$ws.Range("K45").Value = 70000
$ws.Range("M45").Value = -69509
$ws.Range("H62").Value = 205599.2
$ws.Range("I62").Value = 5997
$ws.Range("K62").Value = 5997
$ws.Range("M62").Value = -5373
$ws.Range("H65").Value = 205599.2
$ws.Range("I65").Value = 5997
$ws.Range("K65").Value = 29985
$ws.Range("M65").Value = -26865
$ws.Range("H96").Value = 3319.75
$ws.Range("I96").Value = 2292.375
$ws.Range("K96").Value = 2292.375
$ws.Range("M96").Value = -919.375
$ws.Range("H122").Value = 2062.3333
$ws.Range("I122").Value = 1794.8
$ws.Range("K122").Value = 5384.4
$ws.Range("M122").Value = -2934.4
$ws.Range("H126").Value = 56684.633
$ws.Range("I126").Value = 70373.53
$ws.Range("J126").Value = 5351.25
$ws.Range("K126").Value = 211120.59
$ws.Range("L126").Value = 16053.75
$ws.Range("M126").Value = -208650.59
$ws.Range("N126").Value = -20993.75
$ws.Range("H132").Value = 61602.555
$ws.Range("I132").Value = 62285.117
$ws.Range("K132").Value = 186855.351
$ws.Range("M132").Value = -184325.351
